$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75, shifting existing rows 75-184 down to 76-185
$ws.Range("A75").EntireRow.Insert()

# Fill in the new row 75 with the new record's data.
# Columns A,B,C,E,F,G,H,I,R follow the same boilerplate values as the rest of the sheet.
$ws.Range("A75").Value = 10
$ws.Range("B75").Value = "Vega Modelo de Temuco"
$ws.Range("C75").Value = "La Araucanía"
$ws.Range("D75").Value = 44799
$ws.Range("E75").Value = 9
$ws.Range("F75").Value = 100112012
$ws.Range("G75").Value = "Espinaca"
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 40
$ws.Range("K75").Value = 12000
$ws.Range("L75").Value = 13000
$ws.Range("M75").Value = 12500
$ws.Range("N75").Value = "$/docena de atados"
$ws.Range("O75").Value = "Región de La Araucanía"
$ws.Range("P75").Value = 4167
$ws.Range("Q75").Value = 3
$ws.Range("R75").Value = "Hortaliza"
